{"js": "// Apply the scheduled text replacements (date header + multiplication\n// problems) by locating each unique existing string in the document body\n// and replacing it with its new value.\nconst replacements = [\n  [\"2026-01-24 Saturday\", \"2026-01-25 Sunday\"],\n  [\"41\u00d745=\", \"12\u00d794=\"],\n  [\"23\u00d764=\", \"24\u00d728=\"],\n  [\"45\u00d758=\", \"47\u00d787=\"],\n  [\"81\u00d717=\", \"54\u00d717=\"],\n  [\"63\u00d758=\", \"35\u00d753=\"],\n  [\"47\u00d781=\", \"99\u00d764=\"],\n  [\"52\u00d796=\", \"80\u00d759=\"],\n  [\"39\u00d765=\", \"76\u00d755=\"],\n  [\"13\u00d732=\", \"63\u00d794=\"],\n  [\"68\u00d776=\", \"42\u00d727=\"],\n  [\"34\u00d759=\", \"12\u00d771=\"],\n  [\"42\u00d752=\", \"33\u00d742=\"],\n  [\"16\u00d794=\", \"53\u00d767=\"],\n  [\"37\u00d779=\", \"66\u00d755=\"],\n  [\"13\u00d798=\", \"31\u00d738=\"],\n  [\"25\u00d768=\", \"94\u00d779=\"],\n  [\"71\u00d754=\", \"57\u00d734=\"],\n  [\"80\u00d781=\", \"78\u00d739=\"],\n  [\"27\u00d759=\", \"19\u00d742=\"],\n  [\"58\u00d770=\", \"73\u00d744=\"],\n  [\"35\u00d714=\", \"15\u00d722=\"],\n  [\"36\u00d737=\", \"29\u00d759=\"],\n  [\"69\u00d739=\", \"44\u00d736=\"],\n  [\"96\u00d769=\", \"37\u00d741=\"],\n  [\"63\u00d756=\", \"93\u00d729=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the scheduled text replacements (date header + multiplication\n# problems) using Find/Replace across the whole document story.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-24 Saturday\", \"2026-01-25 Sunday\"),\n    @(\"41\u00d745=\", \"12\u00d794=\"),\n    @(\"23\u00d764=\", \"24\u00d728=\"),\n    @(\"45\u00d758=\", \"47\u00d787=\"),\n    @(\"81\u00d717=\", \"54\u00d717=\"),\n    @(\"63\u00d758=\", \"35\u00d753=\"),\n    @(\"47\u00d781=\", \"99\u00d764=\"),\n    @(\"52\u00d796=\", \"80\u00d759=\"),\n    @(\"39\u00d765=\", \"76\u00d755=\"),\n    @(\"13\u00d732=\", \"63\u00d794=\"),\n    @(\"68\u00d776=\", \"42\u00d727=\"),\n    @(\"34\u00d759=\", \"12\u00d771=\"),\n    @(\"42\u00d752=\", \"33\u00d742=\"),\n    @(\"16\u00d794=\", \"53\u00d767=\"),\n    @(\"37\u00d779=\", \"66\u00d755=\"),\n    @(\"13\u00d798=\", \"31\u00d738=\"),\n    @(\"25\u00d768=\", \"94\u00d779=\"),\n    @(\"71\u00d754=\", \"57\u00d734=\"),\n    @(\"80\u00d781=\", \"78\u00d739=\"),\n    @(\"27\u00d759=\", \"19\u00d742=\"),\n    @(\"58\u00d770=\", \"73\u00d744=\"),\n    @(\"35\u00d714=\", \"15\u00d722=\"),\n    @(\"36\u00d737=\", \"29\u00d759=\"),\n    @(\"69\u00d739=\", \"44\u00d736=\"),\n    @(\"96\u00d769=\", \"37\u00d741=\"),\n    @(\"63\u00d756=\", \"93\u00d729=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
